$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the test case text in F5 (shared string content changes in place,
#    since F5 is the only cell referencing that shared string).
$nl = "`r`n"
$f5Text = '1.Buscar la entidad "HUAWEI_TEST".' + $nl +
          '2.Seleccionar el registro del resultado de búsqueda (card dinámico tipo device).' + $nl +
          '3.Clic en el botón "Editar".' + $nl +
          '4.Clic en la flecha "Siguiente" en el modal de edición.' + $nl +
          '5.Editar campos "Nombre" y "Descripción".' + $nl +
          '6.Clic en el botón "Siguiente" dentro del modal de edición' + $nl +
          '7.Clic en el botón "Editar" y esperar finalización del progreso'
$ws.Range("F5").Value = $f5Text

# 2) Fill in F6 (previously empty) with the new "Eliminar" test case text,
#    which becomes a newly appended shared string.
$f6Text = '1.Buscar la entidad "HUAWEI_TEST_EDIT".' + $nl +
          '2.Seleccionar el registro del resultado de búsqueda (card dinámico tipo device)' + $nl +
          '3.Clic en el botón "Eliminar".' + $nl +
          '4.Clic en el checkbox "Eliminar todas las dependencias"' + $nl +
          '5.Clic en el botón "Eliminar" del modal de confirmación'
$ws.Range("F6").Value = $f6Text

# 3) Make the borders of J6:L6 match the rest of row 6 (full box border, same
#    style as A6) instead of their previous mismatched border styles.
$ws.Range("A6").Copy() | Out-Null
$ws.Range("J6:L6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 4) Update the selected cell shown when the workbook is opened.
$ws.Range("F11").Select() | Out-Null

Write-Host "done"
